# Slide 4, Shape 2 ("TextShape 2") contains the bullet describing the
# "Our solution provides..." sentence that needs to be split into three
# runs so the middle clause "both, probabilistic  database and " gets its
# own run (matching the customized-animation-ready split in the diff).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$oldTail = "both bloom filter and blockchain."
$newTail = "both, probabilistic  database and blockchain."

$fullText = $tr.Text
$startChar = $fullText.IndexOf($oldTail) + 1

# Replace "both bloom filter and blockchain." with the new wording; the
# runtime keeps the surrounding run's formatting for the replaced chunk,
# giving us: [... using combination of ][both, probabilistic  database and blockchain.]
$tailRange = $tr.Characters($startChar, $oldTail.Length)
$tailRange.Text = $newTail

# Now split the tail into two runs so "blockchain." becomes its own run,
# mirroring: [both, probabilistic  database and ][blockchain.]
$fullText2 = $tr.Text
$finalWord = "blockchain."
$finalStart = $fullText2.IndexOf($newTail) + 1 + ($newTail.Length - $finalWord.Length)
$finalRange = $tr.Characters($finalStart, $finalWord.Length)
$finalRange.Text = $finalWord
